# Estrategia.docx: append a new "Migración de la tabla Maestra" section
# (a Title-styled heading paragraph + a body paragraph with a leading tab
# and italic/bold emphasis runs) right after the existing final paragraph,
# before the section properties.
#
# We build the exact target OOXML for the two new <w:p> paragraphs and
# insert it with Range.InsertXML on a range collapsed to the very end of
# the document's story (Content.End). Unlike typing text + toggling
# Bold/Italic/Font.Size, this avoids silently inheriting run/paragraph
# formatting (sz/szCs, indentation, justification, etc.) from the
# preceding paragraph, so the new paragraphs end up with only the
# formatting actually present in the target markup.

$d = $word.ActiveDocument

$newParagraphsXml = '<w:p>' +
    '<w:pPr><w:pStyle w:val="Ttulo"/><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">Migración de la tabla Maestra </w:t></w:r>' +
  '</w:p>' +
  '<w:p>' +
    '<w:pPr><w:rPr><w:u w:val="single"/><w:lang w:val="es-ES"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:tab/><w:t xml:space="preserve">Al migrar los datos de los clientes de la tabla </w:t></w:r>' +
    '<w:r><w:rPr><w:i/><w:lang w:val="es-ES"/></w:rPr><w:t>Maestra</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> hacia la tabla </w:t></w:r>' +
    '<w:r><w:rPr><w:i/><w:lang w:val="es-ES"/></w:rPr><w:t>Clientes</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">, existe la nulidad en el campo </w:t></w:r>' +
    '<w:r><w:rPr><w:i/><w:lang w:val="es-ES"/></w:rPr><w:t>Provincia</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>, por lo cual se carga por defecto con el atributo “</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:lang w:val="es-ES"/></w:rPr><w:t>Migrada</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>”. De esta forma, dejamos en claro, que este cliente fue introducido al sistema a través de la migración y no contaba con una provincia.</w:t></w:r>' +
  '</w:p>'

$openXml = '<?xml version="1.0" standalone="yes"?>' +
  '<?mso-application progid="Word.Document"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' + $newParagraphsXml + '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

# Collapse a range to the end of the main document story (before sectPr)
# and insert the new paragraphs there, after the current last paragraph.
$insertionPoint = $d.Range($d.Content.End, $d.Content.End)
[void]$insertionPoint.InsertXML($openXml)
